# Update the "Metadata" sheet (sheet1) and "Elements" sheet (sheet2) to
# reflect the new FHIR StructureDefinition publication metadata:
#  - Version bump 5.0.0 -> 6.0.0
#  - Date refresh
#  - Publisher value filled in ("Alvearie Team")
#  - "Contact" row repurposed into "Jurisdiction" / "United States of America"
#  - Duplicate "Contact" row removed (rows shift up)
#  - Root extension's Short/Definition text corrected on the Elements sheet

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: refreshed publication timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now populated
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 "Contact" / "No display for ContactDetail" becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row - remove it entirely,
# shifting every following row up by one (dimension becomes A1:B20).
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item("Elements")

# Root Extension element's Short & Definition text updated to describe this specific extension
$elements.Range("K2").Value = "Sent To recipient"
$elements.Range("L2").Value = "Records date/time and outcome of communication attempts with the recipient"
